$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 takes the values previously on row 34
$ws.Range("D2").Value = 44235
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 10533
$ws.Range("N2").Value = '$/malla 18 kilos'
$ws.Range("P2").Value = 585
$ws.Range("Q2").Value = 18

# Row 3 takes the values previously on row 9
$ws.Range("D3").Value = 44382
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 1510
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("N3").Value = '$/malla 18 kilos'
$ws.Range("P3").Value = 472
$ws.Range("Q3").Value = 18

# Row 5 takes the values previously on row 14
$ws.Range("D5").Value = 44403
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 1330
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("N5").Value = '$/caja 15 kilos granel'
$ws.Range("P5").Value = 767
$ws.Range("Q5").Value = 15

# Row 6 takes the values previously on row 30
$ws.Range("D6").Value = 44186
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 1800
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11556
$ws.Range("N6").Value = '$/malla 18 kilos'
$ws.Range("P6").Value = 642
$ws.Range("Q6").Value = 18

# Row 7 takes the values previously on row 6
$ws.Range("D7").Value = 44326
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 1600
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("N7").Value = '$/malla 18 kilos'
$ws.Range("P7").Value = 556
$ws.Range("Q7").Value = 18

# Row 8 takes the values previously on row 2
$ws.Range("D8").Value = 44424
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 790
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13506
$ws.Range("N8").Value = '$/malla 18 kilos'
$ws.Range("P8").Value = 750
$ws.Range("Q8").Value = 18

# Row 9 takes the values previously on row 3
$ws.Range("D9").Value = 44424
$ws.Range("I9").Value = 'Segunda'
$ws.Range("J9").Value = 520
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("N9").Value = '$/malla 18 kilos'
$ws.Range("P9").Value = 667
$ws.Range("Q9").Value = 18

# Row 10 takes the values previously on row 32
$ws.Range("D10").Value = 44179
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 1500
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 10600
$ws.Range("N10").Value = '$/malla 18 kilos'
$ws.Range("P10").Value = 589
$ws.Range("Q10").Value = 18

# Row 11 takes the values previously on row 25
$ws.Range("D11").Value = 44396
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 1330
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9500
$ws.Range("N11").Value = '$/malla 18 kilos'
$ws.Range("P11").Value = 528
$ws.Range("Q11").Value = 18

# Row 12 takes the values previously on row 18
$ws.Range("D12").Value = 44242
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 1600
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = '$/malla 18 kilos'
$ws.Range("P12").Value = 556
$ws.Range("Q12").Value = 18

# Row 13 takes the values previously on row 19
$ws.Range("D13").Value = 44431
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 1150
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12500
$ws.Range("N13").Value = '$/malla 18 kilos'
$ws.Range("P13").Value = 694
$ws.Range("Q13").Value = 18

# Row 14 takes the values previously on row 43
$ws.Range("D14").Value = 44435
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 1150
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 12500
$ws.Range("N14").Value = '$/malla 18 kilos'
$ws.Range("P14").Value = 694
$ws.Range("Q14").Value = 18

# Row 15 takes the values previously on row 12
$ws.Range("D15").Value = 44445
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 1240
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 9500
$ws.Range("N15").Value = '$/malla 18 kilos'
$ws.Range("P15").Value = 528
$ws.Range("Q15").Value = 18

# Row 16 takes the values previously on row 13
$ws.Range("D16").Value = 44445
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 610
$ws.Range("K16").Value = 8000
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = 8000
$ws.Range("N16").Value = '$/malla 18 kilos'
$ws.Range("P16").Value = 444
$ws.Range("Q16").Value = 18

# Row 17 takes the values previously on row 38
$ws.Range("D17").Value = 44263
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 1600
$ws.Range("K17").Value = 9000
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = 9000
$ws.Range("N17").Value = '$/malla 18 kilos'
$ws.Range("P17").Value = 500
$ws.Range("Q17").Value = 18

# Row 18 takes the values previously on row 11
$ws.Range("D18").Value = 44298
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 1600
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 14000
$ws.Range("N18").Value = '$/malla 18 kilos'
$ws.Range("P18").Value = 778
$ws.Range("Q18").Value = 18

# Row 19 takes the values previously on row 44
$ws.Range("D19").Value = 44319
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 1510
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 11000
$ws.Range("M19").Value = 10500
$ws.Range("N19").Value = '$/malla 18 kilos'
$ws.Range("P19").Value = 583
$ws.Range("Q19").Value = 18

# Row 20 takes the values previously on row 16
$ws.Range("D20").Value = 44333
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 1410
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 11000
$ws.Range("M20").Value = 10500
$ws.Range("N20").Value = '$/malla 18 kilos'
$ws.Range("P20").Value = 583
$ws.Range("Q20").Value = 18

# Row 21 takes the values previously on row 24
$ws.Range("D21").Value = 44284
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 1600
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 12000
$ws.Range("N21").Value = '$/malla 18 kilos'
$ws.Range("P21").Value = 667
$ws.Range("Q21").Value = 18

# Row 22 takes the values previously on row 39
$ws.Range("D22").Value = 44417
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 880
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 14500
$ws.Range("N22").Value = '$/malla 18 kilos'
$ws.Range("P22").Value = 806
$ws.Range("Q22").Value = 18

# Row 23 takes the values previously on row 40
$ws.Range("D23").Value = 44417
$ws.Range("I23").Value = 'Segunda'
$ws.Range("J23").Value = 340
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 12000
$ws.Range("N23").Value = '$/malla 18 kilos'
$ws.Range("P23").Value = 667
$ws.Range("Q23").Value = 18

# Row 24 takes the values previously on row 27
$ws.Range("D24").Value = 44340
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 1420
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 11000
$ws.Range("M24").Value = 10500
$ws.Range("N24").Value = '$/malla 18 kilos'
$ws.Range("P24").Value = 583
$ws.Range("Q24").Value = 18

# Row 25 takes the values previously on row 28
$ws.Range("D25").Value = 44340
$ws.Range("I25").Value = 'Segunda'
$ws.Range("J25").Value = 970
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 8000
$ws.Range("N25").Value = '$/malla 18 kilos'
$ws.Range("P25").Value = 444
$ws.Range("Q25").Value = 18

# Row 26 takes the values previously on row 35
$ws.Range("D26").Value = 44452
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 970
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14495
$ws.Range("N26").Value = '$/malla 18 kilos'
$ws.Range("P26").Value = 805
$ws.Range("Q26").Value = 18

# Row 27 takes the values previously on row 36
$ws.Range("D27").Value = 44452
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 340
$ws.Range("K27").Value = 11000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 11500
$ws.Range("N27").Value = '$/malla 18 kilos'
$ws.Range("P27").Value = 639
$ws.Range("Q27").Value = 18

# Row 28 takes the values previously on row 10
$ws.Range("D28").Value = 44193
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 1800
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10556
$ws.Range("N28").Value = '$/malla 18 kilos'
$ws.Range("P28").Value = 586
$ws.Range("Q28").Value = 18

# Row 30 takes the values previously on row 5
$ws.Range("D30").Value = 44214
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 1900
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = 12526
$ws.Range("N30").Value = '$/malla 18 kilos'
$ws.Range("P30").Value = 696
$ws.Range("Q30").Value = 18

# Row 31 takes the values previously on row 15
$ws.Range("D31").Value = 44354
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 700
$ws.Range("K31").Value = 14000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 14500
$ws.Range("N31").Value = '$/malla 18 kilos'
$ws.Range("P31").Value = 806
$ws.Range("Q31").Value = 18

# Row 32 takes the values previously on row 17
$ws.Range("D32").Value = 44165
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 1600
$ws.Range("K32").Value = 11000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = 11625
$ws.Range("N32").Value = '$/malla 18 kilos'
$ws.Range("P32").Value = 646
$ws.Range("Q32").Value = 18

# Row 33 takes the values previously on row 21
$ws.Range("D33").Value = 44172
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 1600
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 12000
$ws.Range("N33").Value = '$/malla 18 kilos'
$ws.Range("P33").Value = 667
$ws.Range("Q33").Value = 18

# Row 34 takes the values previously on row 8
$ws.Range("D34").Value = 44389
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 1420
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = 8500
$ws.Range("N34").Value = '$/malla 18 kilos'
$ws.Range("P34").Value = 472
$ws.Range("Q34").Value = 18

# Row 35 takes the values previously on row 33
$ws.Range("D35").Value = 44249
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 1600
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 10000
$ws.Range("N35").Value = '$/malla 18 kilos'
$ws.Range("P35").Value = 556
$ws.Range("Q35").Value = 18

# Row 36 takes the values previously on row 7
$ws.Range("D36").Value = 44270
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 16000
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = 10500
$ws.Range("N36").Value = '$/malla 18 kilos'
$ws.Range("P36").Value = 583
$ws.Range("Q36").Value = 18

# Row 37 takes the values previously on row 20
$ws.Range("D37").Value = 44312
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 1510
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 12000
$ws.Range("M37").Value = 12000
$ws.Range("N37").Value = '$/malla 18 kilos'
$ws.Range("P37").Value = 667
$ws.Range("Q37").Value = 18

# Row 38 takes the values previously on row 26
$ws.Range("D38").Value = 44221
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 1800
$ws.Range("K38").Value = 11000
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 11444
$ws.Range("N38").Value = '$/malla 18 kilos'
$ws.Range("P38").Value = 636
$ws.Range("Q38").Value = 18

# Row 39 takes the values previously on row 31
$ws.Range("D39").Value = 44277
$ws.Range("I39").Value = 'Primera'
$ws.Range("J39").Value = 1600
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 15000
$ws.Range("N39").Value = '$/malla 18 kilos'
$ws.Range("P39").Value = 833
$ws.Range("Q39").Value = 18

# Row 40 takes the values previously on row 37
$ws.Range("D40").Value = 44291
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 1600
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 12000
$ws.Range("N40").Value = '$/malla 18 kilos'
$ws.Range("P40").Value = 667
$ws.Range("Q40").Value = 18

# Row 43 takes the values previously on row 22
$ws.Range("D43").Value = 44410
$ws.Range("I43").Value = 'Primera'
$ws.Range("J43").Value = 970
$ws.Range("K43").Value = 14000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 14505
$ws.Range("N43").Value = '$/malla 18 kilos'
$ws.Range("P43").Value = 806
$ws.Range("Q43").Value = 18

# Row 44 takes the values previously on row 23
$ws.Range("D44").Value = 44410
$ws.Range("I44").Value = 'Segunda'
$ws.Range("J44").Value = 340
$ws.Range("K44").Value = 12000
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = 12000
$ws.Range("N44").Value = '$/malla 18 kilos'
$ws.Range("P44").Value = 667
$ws.Range("Q44").Value = 18

